$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 00:52"

# Row 4
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 736790
$ws.Range("C4").Value = 27055
$ws.Range("D4").Value = 67438
$ws.Range("E4").Value = 630432
$ws.Range("F4").Value = 13536
$ws.Range("G4").Value = 1766
$ws.Range("H4").Value = 38920

# Row 5
$ws.Range("A5").Value = "España"
$ws.Range("B5").Value = 191726
$ws.Range("C5").Value = 887
$ws.Range("D5").Value = 74797
$ws.Range("E5").Value = 96290
$ws.Range("F5").Value = 7371
$ws.Range("G5").Value = 637
$ws.Range("H5").Value = 20639

# Row 51
$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 3621
$ws.Range("C51").Value = 182
$ws.Range("D51").Value = 691
$ws.Range("E51").Value = 2764
$ws.Range("F51").Value = 98
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 166

# Row 52
$ws.Range("A52").Value = "Luxemburgo"
$ws.Range("B52").Value = 3537
$ws.Range("C52").Value = 57
$ws.Range("D52").Value = 601
$ws.Range("E52").Value = 2864
$ws.Range("F52").Value = 32
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 72

# Row 85
$ws.Range("A85").Value = "Tunez"
$ws.Range("B85").Value = 866
$ws.Range("C85").Value = 2
$ws.Range("D85").Value = 43
$ws.Range("E85").Value = 786
$ws.Range("F85").Value = 33
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 37

# Row 98
$ws.Range("A98").Value = "Nigeria"
$ws.Range("B98").Value = 542
$ws.Range("C98").Value = 49
$ws.Range("D98").Value = 166
$ws.Range("E98").Value = 357
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 19

# Row 99
$ws.Range("A99").Value = "Guinea"
$ws.Range("B99").Value = 518
$ws.Range("C99").Value = 41
$ws.Range("D99").Value = 65
$ws.Range("E99").Value = 450
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 3

# Row 100
$ws.Range("A100").Value = "Uruguay"
$ws.Range("B100").Value = 508
$ws.Range("C100").Value = 6
$ws.Range("D100").Value = 294
$ws.Range("E100").Value = 205
$ws.Range("F100").Value = 12
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 9

# Row 101
$ws.Range("A101").Value = "Kirguistan"
$ws.Range("B101").Value = 506
$ws.Range("C101").Value = 17
$ws.Range("D101").Value = 130
$ws.Range("E101").Value = 371
$ws.Range("F101").Value = 5
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 5

# Row 102
$ws.Range("A102").Value = "Bolivia"
$ws.Range("B102").Value = 493
$ws.Range("C102").Value = 28
$ws.Range("D102").Value = 31
$ws.Range("E102").Value = 431
$ws.Range("F102").Value = 3
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 31

# Row 163
$ws.Range("A163").Value = "Puerto Rico"
$ws.Range("B163").Value = 39
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 1
$ws.Range("E163").Value = 36
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 2

# Row 164
$ws.Range("A164").Value = "Eritrea"
$ws.Range("B164").Value = 39
$ws.Range("C164").Value = 4
$ws.Range("D164").Value = 3
$ws.Range("E164").Value = 36
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 0

# Row 167
$ws.Range("A167").Value = "Mozambique"
$ws.Range("B167").Value = 35
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 4
$ws.Range("E167").Value = 31
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

# Row 168
$ws.Range("A168").Value = "Maldivas"
$ws.Range("B168").Value = 35
$ws.Range("C168").Value = 6
$ws.Range("D168").Value = 16
$ws.Range("E168").Value = 19
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

# Row 169
$ws.Range("A169").Value = "Benin"
$ws.Range("B169").Value = 35
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 18
$ws.Range("E169").Value = 16
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0
